# Updated on 12.24pm 24th Feb 2015 from SR
#
# Marks a set of task rows in column C with "o" (done marker) and moves the
# sheet's active selection/scroll position to reflect where the editor was
# last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Mark these rows as done ("o") in column C.
$doneRows = @(27, 28, 30, 31, 32, 33, 35, 39, 41)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 3).Value = "o"
}

# Reflect the editor's final scroll/selection position.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("K28").Select()
